$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Difficulty" / "Bombs on?" columns aren't used by the reader - clear them out.
$ws.Range("C1:D1").ClearContents()

# Add two new high-score rows read from the game.
$ws.Range("A2").Value = "Jack"
$ws.Range("B2").Value = 5
$ws.Range("A3").Value = "Jack Better"
$ws.Range("B3").Value = 10

# Sort the new rows by score, highest first.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B3"), 0, 2, 0, 0) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:B3"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Leave the selection where the author left it.
$ws.Range("C1").Select() | Out-Null
